$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.279.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.25%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.620.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.71%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'519.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.04%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'150.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.89%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.13%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'6.39"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.78%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.106"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.02%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.344"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.47%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.93%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.079.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.54%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'60.283.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.35%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'21.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.57%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000139"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.615.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.12%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.87%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'346.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.96%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.63%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.32%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.36%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'60.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.16%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.421"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.164"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.94%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0₃0833"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.74%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.69%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.00%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.63%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.89%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'19.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.03%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'149.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.90%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'1.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.99%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.890"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.77%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.878"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +4.18%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'36.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.20%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.94%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Filecoin"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'3.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.74%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Bittensor"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'289.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.629"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.17%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.100"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.17%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.06%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0552"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.61%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'19.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.15%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0235"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.35%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'4.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.97%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'10.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.50%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'18.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.46%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.959.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.43%  "
$ws.Range("E51").Style = "Normal"
